$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the last existing data row (162) down into the two
# new rows so the date column keeps its "yyyy-mm-dd hh:mm:ss" style (s="1")
# and the other columns keep the default (unstyled) look.
$ws.Range("A162:H162").Copy()
$ws.Range("A163:H163").PasteSpecial(-4122)
$ws.Range("A164:H164").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 163 ---
$ws.Range("A163").Value = 45477.2916666667
$ws.Range("B163").Value = 0
$ws.Range("C163").Value = 7.55000019073486
$ws.Range("D163").Value = 7.55000019073486
$ws.Range("E163").Value = 7.55000019073486
$ws.Range("F163").Value = 7.55000019073486

# adj_close is stored as text in this workbook, so force text entry and
# then drop the temporary "Text" number format again (real cells in the
# existing data use General format, just with a string value).
$ws.Range("G163").NumberFormat = "@"
$ws.Range("G163").Value = "7.55000019073486"
$ws.Range("G163").ClearFormats()

$ws.Range("H163").Value = "VARV.MI"

# --- Row 164 ---
$ws.Range("A164").Value = 45478.6298842593
$ws.Range("B164").Value = 600
$ws.Range("C164").Value = 7.40000009536743
$ws.Range("D164").Value = 7.40000009536743
$ws.Range("E164").Value = 7.40000009536743
$ws.Range("F164").Value = 7.40000009536743

$ws.Range("G164").NumberFormat = "@"
$ws.Range("G164").Value = "7.40000009536743"
$ws.Range("G164").ClearFormats()

$ws.Range("H164").Value = "VARV.MI"
